$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell -> new text value map, derived from the authoritative XML diff.
# Every one of these cells is stored as inline text (t="inlineStr") in the
# workbook, even the ones that look numeric/percentage, so each write forces
# Text formatting first (NumberFormat "@") and then restores the default
# "Normal" style afterwards so no stray style/number-format is left behind.
$changes = @{
    'D2' = '304.23'
    'E2' = '2.07%'
    'D3' = '31.90'
    'E3' = '0.11%'
    'D4' = '5.217'
    'E4' = '2.04%'
    'D5' = '0.07822'
    'E5' = '3.75%'
    'D6' = '2.314'
    'E6' = '32.85%'
    'D7' = '7.995'
    'E7' = '3.11%'
    'D8' = '3.871'
    'E8' = '2.07%'
    'D9' = '0.9135'
    'E9' = '-1.68%'
    'D10' = '0.1741'
    'E10' = '2.44%'
    'D11' = '0.07388'
    'E11' = '-1.49%'
    'D12' = '0.08163'
    'E12' = '2.59%'
    'D13' = '0.03045'
    'E13' = '-0.27%'
    'E14' = '0.76%'
    'D15' = '0.001517'
    'E15' = '1.53%'
    'D16' = '0.006180'
    'E16' = '-4.44%'
    'D17' = '3.500'
    'E17' = '1.12%'
    'E18' = '0.86%'
    'D19' = '0.3279'
    'E19' = '-0.06%'
    'D20' = '0.1329'
    'E20' = '0.31%'
    'D21' = '4.658'
    'E21' = '1.96%'
    'D22' = '0.04634'
    'E22' = '-0.42%'
    'D23' = '0.1565'
    'E23' = '0.52%'
    'E24' = '3.51%'
    'D25' = '0.004531'
    'E25' = '2.42%'
    'D26' = '0.0001349'
    'E26' = '-3.62%'
    'D27' = '0.0002740'
    'E27' = '47.49%'
    'D39' = '0.01790'
    'E39' = '6.41%'
    'D40' = '0.04591'
    'E40' = '1.25%'
    'E41' = '2.76%'
    'D42' = '0.1364'
    'E42' = '2.86%'
    'D43' = '0.002239'
    'E43' = '8.78%'
    'D44' = '0.01088'
    'E44' = '-6.51%'
    'D45' = '0.00006378'
    'E45' = '7.41%'
    'B46' = 'Kangarootoken'
    'C46' = 'https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar'
    'D46' = '0.00000000749'
    'E46' = '-0.10%'
    'B47' = 'BOLO'
    'C47' = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
    'D47' = '0.8206'
    'E47' = '-57.48%'
    'B48' = 'CoinbaseStockToken'
    'C48' = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
    'D48' = '0.009893'
    'E48' = '-23.70%'
    'B49' = 'CryptobidCoin'
    'C49' = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
    'D49' = '0.00002098'
    'E49' = '-0.10%'
    'B50' = 'SpecialPowerGold'
    'C50' = 'https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg'
    'D50' = '0.0001998'
    'E50' = '-0.02%'
    'B51' = 'DigiFinexToken'
    'C51' = 'https://coinranking.com/coin/rY6dWXQL4+digifinextoken-dft'
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.Style = "Normal"
}
